# Refactor import service fixture: drop the second review sheet, swap the
# numeric b24StationId values for formatted service labels on the remaining
# sheet, and append a third review row.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Delete the now-unused "reviews_test2" sheet first so the strings that
# only it referenced ("Отзыв 5" / "Бяка") are freed up before we add the
# new ones.
$wb.Worksheets("reviews_test2").Delete()

$ws = $wb.Worksheets("review_test")

# Row 2: b24StationId 138 becomes a descriptive string value instead of a
# bare number.
$ws.Range("A2").Value = "[138] Service Архангельск"

# Row 3: b24StationId 178 becomes a descriptive string value instead of a
# bare number.
$ws.Range("A3").Value = "[178] Service — Колпино"

# New row 4: another review for station 178, re-using the label above.
$ws.Range("A4").Value = "[178] Service — Колпино"
$ws.Range("B4").Value = "Отзыв 3"
$ws.Range("C4").Value = "Бука"
# F4 stays empty but keeps the date style used by the rest of column F.
$ws.Range("F4").NumberFormat = "yyyy\-mm\-dd"

# Column width tweaks.
$ws.Columns.Item(1).ColumnWidth = 11.45
$ws.Columns.Item(5).ColumnWidth = 23.08
$ws.StandardWidth = 8.58984375

# Move the selection to A5, below the new data.
$ws.Range("A5").Select()
